$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1461.5454
$ws.Range("I2").Value = 709.625
$ws.Range("J2").Value = 3466.6667
$ws.Range("K2").Value = 709.625
$ws.Range("L2").Value = 3466.6667
$ws.Range("M2").Value = -596.625
$ws.Range("N2").Value = -3692.6667
$ws.Range("H62").Value = 5487.647
$ws.Range("I62").Value = 5153
$ws.Range("K62").Value = 5153
$ws.Range("M62").Value = -4529
$ws.Range("H65").Value = 5487.647
$ws.Range("I65").Value = 5153
$ws.Range("K65").Value = 25765
$ws.Range("M65").Value = -22645
$ws.Range("H100").Value = 3449.75
$ws.Range("I100").Value = 3449.75
$ws.Range("K100").Value = 3449.75
$ws.Range("M100").Value = -2908.75
$ws.Range("H107").Value = 1427.7715
$ws.Range("I107").Value = 634.4231
$ws.Range("J107").Value = 3719.6667
$ws.Range("K107").Value = 634.4231
$ws.Range("L107").Value = 3719.6667
$ws.Range("M107").Value = 1285.5769
$ws.Range("N107").Value = -7559.6667
$ws.Range("H137").Value = 2922.3076
$ws.Range("J137").Value = 3294
$ws.Range("L137").Value = 9882
$ws.Range("N137").Value = -14982

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1841.8422
$ws.Range("I2").Value = 1213
$ws.Range("K2").Value = 1213
$ws.Range("M2").Value = -1100
$ws.Range("H32").Value = 7076.735
$ws.Range("I32").Value = 7076.735
$ws.Range("K32").Value = 7076.735
$ws.Range("M32").Value = -6789.735
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H97").Value = 927.2353000000001
$ws.Range("I97").Value = 664.2143
$ws.Range("J97").Value = 2154.6667
$ws.Range("K97").Value = 664.2143
$ws.Range("L97").Value = 2154.6667
$ws.Range("M97").Value = -168.2143
$ws.Range("N97").Value = -3146.6667
$ws.Range("H102").Value = 816.2222
$ws.Range("I102").Value = 793.25
$ws.Range("K102").Value = 793.25
$ws.Range("M102").Value = 828.75
$ws.Range("H110").Value = 2246.9333
$ws.Range("I110").Value = 1120.4
$ws.Range("K110").Value = 1120.4
$ws.Range("M110").Value = 924.5999999999999
$ws.Range("H116").Value = 1841.8422
$ws.Range("I116").Value = 1213
$ws.Range("K116").Value = 1213
$ws.Range("M116").Value = 1081
$ws.Range("H132").Value = 3776
$ws.Range("I132").Value = 4154.857
$ws.Range("K132").Value = 12464.571
$ws.Range("M132").Value = -9934.571
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1841.8422
$ws.Range("I3").Value = 1213
$ws.Range("K3").Value = 1213
$ws.Range("M3").Value = -1099
$ws.Range("H20").Value = 4666.5
$ws.Range("I20").Value = 4666.5
$ws.Range("K20").Value = 4666.5
$ws.Range("M20").Value = -4419.5
$ws.Range("H99").Value = 4989.6665
$ws.Range("I99").Value = 4989.6665
$ws.Range("K99").Value = 4989.6665
$ws.Range("M99").Value = -3491.6665
$ws.Range("H107").Value = 1172.25
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("H108").Value = 79940
$ws.Range("J108").Value = 79940
$ws.Range("L108").Value = 79940
$ws.Range("N108").Value = -87620
$ws.Range("H115").Value = 40145
$ws.Range("J115").Value = 40145
$ws.Range("L115").Value = 40145
$ws.Range("N115").Value = -42495
$ws.Range("H134").Value = 3087.9546
$ws.Range("I134").Value = 2513.4707
$ws.Range("J134").Value = 5041.2
$ws.Range("K134").Value = 7540.4121
$ws.Range("L134").Value = 15123.6
$ws.Range("M134").Value = -5005.4121
$ws.Range("N134").Value = -20193.6
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1479.75
$ws.Range("J32").Value = 2910
$ws.Range("L32").Value = 8730
$ws.Range("N32").Value = -9296
$ws.Range("H46").Value = 50
$ws.Range("I46").Value = 50
$ws.Range("K46").Value = 150
$ws.Range("M46").Value = -59
$ws.Range("H97").Value = 834.3077
$ws.Range("J97").Value = 535.1429000000001
$ws.Range("L97").Value = 1605.4287
$ws.Range("N97").Value = -2597.4287
$ws.Range("H132").Value = 1200
$ws.Range("I132").Value = 1200
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10800
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8270
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 19125
$ws.Range("I35").Value = 19125
$ws.Range("K35").Value = 19125
$ws.Range("M35").Value = -18827
$ws.Range("H102").Value = 2493.375
$ws.Range("I102").Value = 2663.8572
$ws.Range("K102").Value = 2663.8572
$ws.Range("M102").Value = -1041.8572
$ws.Range("H126").Value = 4600
$ws.Range("I126").Value = 4600
$ws.Range("J126").Value = 4600
$ws.Range("K126").Value = 13800
$ws.Range("L126").Value = 13800
$ws.Range("M126").Value = -11330
$ws.Range("N126").Value = -18740
$ws.Range("H132").Value = 2392.2
$ws.Range("I132").Value = 2392.2
$ws.Range("K132").Value = 7176.599999999999
$ws.Range("M132").Value = -4646.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 599.5
$ws.Range("I7").Value = 599.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 599.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -487.5
$ws.Range("N7").ClearContents()
$ws.Range("H16").Value = 3569.8
$ws.Range("I16").Value = 3633.111
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 3633.111
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -3463.111
$ws.Range("N16").Value = -3340
$ws.Range("H40").Value = 3844.111
$ws.Range("I40").Value = 3712.25
$ws.Range("K40").Value = 3712.25
$ws.Range("M40").Value = -3576.25
$ws.Range("H82").Value = 3155.2
$ws.Range("I82").Value = 2638
$ws.Range("J82").Value = 3500
$ws.Range("K82").Value = 2638
$ws.Range("L82").Value = 3500
$ws.Range("M82").Value = -2277
$ws.Range("N82").Value = -4222
$ws.Range("H85").Value = 3155.2
$ws.Range("I85").Value = 2638
$ws.Range("J85").Value = 3500
$ws.Range("K85").Value = 2638
$ws.Range("L85").Value = 3500
$ws.Range("M85").Value = -1390
$ws.Range("N85").Value = -5996
$ws.Range("H126").Value = 599.5
$ws.Range("I126").Value = 599.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 1798.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 671.5
$ws.Range("N126").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5082.1665
$ws.Range("I126").Value = 5082.1665
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15246.4995
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12776.4995
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 1861.6111
$ws.Range("I136").Value = 1735.8235
$ws.Range("K136").Value = 5207.470499999999
$ws.Range("M136").Value = -2657.470499999999
